$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 24,13

$data[0,0] = 0.7710948365625825
$data[0,1] = 0.05543888078098291
$data[0,2] = 0.07738518413988515
$data[0,3] = 0.07458241688745293
$data[0,4] = 0
$data[0,5] = 1.615497970026766
$data[0,6] = 1.440838659402161
$data[0,7] = 1.416232152607925
$data[0,8] = 0
$data[0,9] = 0.5687447465496405
$data[0,10] = 0.2201298687743147
$data[0,11] = 0.2087863260559963
$data[0,12] = 2.639224997616026

$data[1,0] = 0.7348842129705417
$data[1,1] = 0.04972933416826208
$data[1,2] = 0.07036235194044593
$data[1,3] = 0.07478131441755131
$data[1,4] = 0
$data[1,5] = 1.609757131643121
$data[1,6] = 1.443228934598679
$data[1,7] = 1.41833671836158
$data[1,8] = 0
$data[1,9] = 0.5288870036945923
$data[1,10] = 0.217449827957509
$data[1,11] = 0.2020892534365508
$data[1,12] = 2.65760901537395

$data[2,0] = 0.7130540906033787
$data[2,1] = 0.04619936719932127
$data[2,2] = 0.06608662543308697
$data[2,3] = 0.07492304629010427
$data[2,4] = 0
$data[2,5] = 1.607037261016785
$data[2,6] = 1.445254469117359
$data[2,7] = 1.42021005889805
$data[2,8] = 0
$data[2,9] = 0.5046943151792789
$data[2,10] = 0.215911052343138
$data[2,11] = 0.1980878122596792
$data[2,12] = 2.669618202713536

$data[3,0] = 0.7042598741784332
$data[3,1] = 0.04475471512880347
$data[3,2] = 0.0643533333743278
$data[3,3] = 0.07498574539312397
$data[3,4] = 0
$data[3,5] = 1.606131225596783
$data[3,6] = 1.44622025541301
$data[3,7] = 1.421119616932458
$data[3,8] = 0
$data[3,9] = 0.4949061942735113
$data[3,10] = 0.2153108813087456
$data[3,11] = 0.1964850665398252
$data[3,12] = 2.674693396369904

$data[4,0] = 0.7028057550494111
$data[4,1] = 0.04451445796252074
$data[4,2] = 0.0640660703131033
$data[4,3] = 0.07499645537145216
$data[4,4] = 0
$data[4,5] = 1.60599299681698
$data[4,6] = 1.446389103470921
$data[4,7] = 1.421279476685655
$data[4,8] = 0
$data[4,9] = 0.4932851527599382
$data[4,10] = 0.215212848845141
$data[4,11] = 0.1962206173061105
$data[4,12] = 2.675547080788135

$data[5,0] = 0.712935076387339
$data[5,1] = 0.04617990913833125
$data[5,2] = 0.06606321284132832
$data[5,3] = 0.07492387184565352
$data[5,4] = 0
$data[5,5] = 1.607024222795189
$data[5,6] = 1.445266925620459
$data[5,7] = 1.42022173370102
$data[5,8] = 0
$data[5,9] = 0.5045620231082069
$data[5,10] = 0.2159028492833812
$data[5,11] = 0.1980660841219368
$data[5,12] = 2.669685914310243

$data[6,0] = 0.7585258717564329
$data[6,1] = 0.05347523375689889
$data[6,2] = 0.07495615176503634
$data[6,3] = 0.07464693403317924
$data[6,4] = 0
$data[6,5] = 1.613351380963266
$data[6,6] = 1.441547078514702
$data[6,7] = 1.4168372069661
$data[6,8] = 0
$data[6,9] = 0.5549437361465266
$data[6,10] = 0.219183652059705
$data[6,11] = 0.2064542594586811
$data[6,12] = 2.645414077866022

$data[7,0] = 0.8511215201824598
$data[7,1] = 0.06759180568457168
$data[7,2] = 0.09268574400189777
$data[7,3] = 0.07425895221900802
$data[7,4] = 0
$data[7,5] = 1.632153551189091
$data[7,6] = 1.438676884950908
$data[7,7] = 1.414811217844736
$data[7,8] = 0
$data[7,9] = 0.6559636486549323
$data[7,10] = 0.2264633132809024
$data[7,11] = 0.2237791389774273
$data[7,12] = 2.603541517433236

$data[8,0] = 0.9210948356803783
$data[8,1] = 0.07785291635715907
$data[8,2] = 0.1058934478882207
$data[8,3] = 0.0740678275651252
$data[8,4] = 0
$data[8,5] = 1.649879630635468
$data[8,6] = 1.43926409276753
$data[8,7] = 1.416135690323394
$data[8,8] = 0
$data[8,9] = 0.7315441646653085
$data[8,10] = 0.232326657349347
$data[8,11] = 0.2370408256643657
$data[8,12] = 2.576268443958973

$data[9,0] = 0.9533494237016669
$data[9,1] = 0.08249817304462681
$data[9,2] = 0.1119425170096662
$data[9,3] = 0.07400114624621423
$data[9,4] = 0
$data[9,5] = 1.658796420599032
$data[9,6] = 1.440116504347827
$data[9,7] = 1.417349611978281
$data[9,8] = 0
$data[9,9] = 0.7662252760060824
$data[9,10] = 0.2351057628463735
$data[9,11] = 0.243189617845303
$data[9,12] = 2.56461944374194

$data[10,0] = 0.9656240932583842
$data[10,1] = 0.0842540383453354
$data[10,2] = 0.1142390672556104
$data[10,3] = 0.07397879845428257
$data[10,4] = 0
$data[10,5] = 1.662295845295461
$data[10,6] = 1.440523419376149
$data[10,7] = 1.417897232625315
$data[10,8] = 0
$data[10,9] = 0.7794011059746708
$data[10,10] = 0.2361741934350761
$data[10,11] = 0.245534644075299
$data[10,12] = 2.560317277170761

$data[11,0] = 0.962977835395435
$data[11,1] = 0.08387602247552195
$data[11,2] = 0.1137442012419143
$data[11,3] = 0.07398348251033404
$data[11,4] = 0
$data[11,5] = 1.661536716209639
$data[11,6] = 1.440432042041977
$data[11,7] = 1.417775381545731
$data[11,8] = 0
$data[11,9] = 0.7765615524028533
$data[11,10] = 0.2359433748634387
$data[11,11] = 0.2450288629453397
$data[11,12] = 2.561238973638744

$data[12,0] = 0.9543580563494061
$data[12,1] = 0.08264269281791314
$data[12,2] = 0.1121313371873498
$data[12,3] = 0.07399924956106751
$data[12,4] = 0
$data[12,5] = 1.659081857793296
$data[12,6] = 1.440148295710173
$data[12,7] = 1.417392902624513
$data[12,8] = 0
$data[12,9] = 0.7673084017584131
$data[12,10] = 0.2351933419571566
$data[12,11] = 0.2433822120127829
$data[12,12] = 2.564263315620664

$data[13,0] = 0.9490860695257766
$data[13,1] = 0.0818868286095551
$data[13,2] = 0.1111441804827251
$data[13,3] = 0.07400928507126991
$data[13,4] = 0
$data[13,5] = 1.657594186806108
$data[13,6] = 1.43998544737164
$data[13,7] = 1.417170075412443
$data[13,8] = 0
$data[13,9] = 0.7616461544936044
$data[13,10] = 0.2347360133509824
$data[13,11] = 0.2423757521894885
$data[13,12] = 2.566130018078297

$data[14,0] = 0.9189954141875205
$data[14,1] = 0.07754888823667727
$data[14,2] = 0.1054989508860302
$data[14,3] = 0.07407259215741568
$data[14,4] = 0
$data[14,5] = 1.649314077270049
$data[14,6] = 1.439220159828636
$data[14,7] = 1.416068662073521
$data[14,8] = 0
$data[14,9] = 0.7292836774911962
$data[14,10] = 0.232147283461714
$data[14,11] = 0.2366413158599556
$data[14,12] = 2.577044988571494

$data[15,0] = 0.9006439776253501
$data[15,1] = 0.07488195866240233
$data[15,2] = 0.1020462609564277
$data[15,3] = 0.07411661101175326
$data[15,4] = 0
$data[15,5] = 1.644453108582724
$data[15,6] = 1.438900559072977
$data[15,7] = 1.415549587893388
$data[15,8] = 0
$data[15,9] = 0.7095068440102636
$data[15,10] = 0.2305877988929694
$data[15,11] = 0.2331530840035683
$data[15,12] = 2.583935098581406

$data[16,0] = 0.8901285861005306
$data[16,1] = 0.07334588333863223
$data[16,2] = 0.1000642053097351
$data[16,3] = 0.07414383730797702
$data[16,4] = 0
$data[16,5] = 1.641737497139019
$data[16,6] = 1.438771825395747
$data[16,7] = 1.415308578211011
$data[16,8] = 0
$data[16,9] = 0.6981599006562078
$data[16,10] = 0.2297013539248098
$data[16,11] = 0.231157666247924
$data[16,12] = 2.587969415313651

$data[17,0] = 0.8865751127104033
$data[17,1] = 0.07282542751464405
$data[17,2] = 0.09939377399345517
$data[17,3] = 0.07415338366038426
$data[17,4] = 0
$data[17,5] = 1.640831822851709
$data[17,6] = 1.438737701733032
$data[17,7] = 1.41523685999757
$data[17,8] = 0
$data[17,9] = 0.6943228637849757
$data[17,10] = 0.2294030284973871
$data[17,11] = 0.2304839297912125
$data[17,12] = 2.589347609621292

$data[18,0] = 0.9025933953917615
$data[18,1] = 0.07516607771880501
$data[18,2] = 0.1024134080255692
$data[18,3] = 0.07411172775580788
$data[18,4] = 0
$data[18,5] = 1.644962256254985
$data[18,6] = 1.438928879559683
$data[18,7] = 1.415598888113749
$data[18,8] = 0
$data[18,9] = 0.7116092076475411
$data[18,10] = 0.2307527192265582
$data[18,11] = 0.2335232824441391
$data[18,12] = 2.583194253759082

$data[19,0] = 0.9568882534735224
$data[19,1] = 0.08300503785443425
$data[19,2] = 0.1126049141253134
$data[19,3] = 0.0739945396951498
$data[19,4] = 0
$data[19,5] = 1.659799574412745
$data[19,6] = 1.440229356048263
$data[19,7] = 1.417502859304896
$data[19,8] = 0
$data[19,9] = 0.7700251149497888
$data[19,10] = 0.2354132096920551
$data[19,11] = 0.2438654229873265
$data[19,12] = 2.563372032425256

$data[20,0] = 0.992725905683983
$data[20,1] = 0.08810966742423432
$data[20,2] = 0.1193000501718728
$data[20,3] = 0.07393486636346935
$data[20,4] = 0
$data[20,5] = 1.670212638771432
$data[20,6] = 1.441569622461401
$data[20,7] = 1.419259774055845
$data[20,8] = 0
$data[20,9] = 0.8084529459416103
$data[20,10] = 0.238552606379983
$data[20,11] = 0.2507214070135575
$data[20,12] = 2.551052784579824

$data[21,0] = 0.9735664994717581
$data[21,1] = 0.0853869133088665
$data[21,2] = 0.115723573015984
$data[21,3] = 0.07396517086595544
$data[21,4] = 0
$data[21,5] = 1.664589423167143
$data[21,6] = 1.440809444076734
$data[21,7] = 1.418275171169974
$data[21,8] = 0
$data[21,9] = 0.7879205022701967
$data[21,10] = 0.2368685091357889
$data[21,11] = 0.2470534041898986
$data[21,12] = 2.557569594435392

$data[22,0] = 0.9017119542817227
$data[22,1] = 0.0750376362889682
$data[22,2] = 0.1022474116770411
$data[22,3] = 0.07411392949700968
$data[22,4] = 0
$data[22,5] = 1.644731824442232
$data[22,6] = 1.438915904505961
$data[22,7] = 1.415576420640427
$data[22,8] = 0
$data[22,9] = 0.7106586573268601
$data[22,10] = 0.2306781272079661
$data[22,11] = 0.2333558845344612
$data[22,12] = 2.583528961990872

$data[23,0] = 0.8257304863251989
$data[23,1] = 0.06379265141977442
$data[23,2] = 0.08785783244113077
$data[23,3] = 0.07434737249262291
$data[23,4] = 0
$data[23,5] = 1.6263810592027
$data[23,6] = 1.438979825926666
$data[23,7] = 1.414865461609772
$data[23,8] = 0
$data[23,9] = 0.6283965485025362
$data[23,10] = 0.2244034813255666
$data[23,11] = 0.2189986296512529
$data[23,12] = 2.614256232130202

$ws.Range("B2:N25").Value = $data

Write-Output "Updated data rows 2-25 (columns B-N)"
